# Add a new slide (5th) after the existing 4 slides, using the
# "Title and Content" layout (same layout used by slides 2-4).
$p = $ppt.ActivePresentation
$s = $p.Slides.Add(5, 2)

# Set the slide title.
$s.Shapes.Item(1).TextFrame.TextRange.Text = "TODO tech that’s new to me, ones I know"

# Replace the (empty) content placeholder with a 6x5 table occupying the
# same area the placeholder would have used.
$s.Shapes.Item(2).Delete()
$tbl = $s.Shapes.AddTable(6, 5, 66, 143.75, 828, 175.2)
$tbl.Name = "Content Placeholder 3"

$table = $tbl.Table

# Row 1: header row, left blank.

# Row 2
$table.Cell(2,1).Shape.TextFrame.TextRange.Text = "C#"
$table.Cell(2,2).Shape.TextFrame.TextRange.Text = "Java"

# Row 3
$table.Cell(3,1).Shape.TextFrame.TextRange.Text = "LINQ"
$table.Cell(3,2).Shape.TextFrame.TextRange.Text = "Hibernate"

# Row 4
$table.Cell(4,1).Shape.TextFrame.TextRange.Text = "Visual Studio"
$table.Cell(4,2).Shape.TextFrame.TextRange.Text = "Eclipse"

# Row 5
$table.Cell(5,1).Shape.TextFrame.TextRange.Text = "REACT"
$table.Cell(5,2).Shape.TextFrame.TextRange.Text = "JSPs"

# Row 6
$table.Cell(6,1).Shape.TextFrame.TextRange.Text = "TODO any else?"
